# Daily attendance processing - reorder "Recorded By" (column G) entries
# so that "System"/"system" entries are moved to the end of the
# comma-separated list, preserving the relative order of the other entries.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2

    if ($val -ne $null -and $val -ne "") {
        $parts = $val -split ", "

        if ($parts.Length -gt 1) {
            $otherParts = @()
            $systemParts = @()

            foreach ($p in $parts) {
                if ($p -eq "System" -or $p -eq "system") {
                    $systemParts += $p
                } else {
                    $otherParts += $p
                }
            }

            $newVal = ($otherParts + $systemParts) -join ", "

            if ($newVal -ne $val) {
                $cell.Value = $newVal
            }
        }
    }
}
